$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number and must be forced back to
# Text (the sheet stores prices/volumes as text, e.g. "1.00", not numeric 1)
# so typing them would otherwise be auto-converted to a Number by Excel.
$forceText = @("D4", "D5", "D6", "D14", "D19", "D20", "D21", "D24", "D26", "D27", "D40", "D41", "D44", "D45", "D47", "D48", "D49", "D50")

$updates = @{
    "D2" = "62.194.34"
    "E2" = "  +0.99%  "
    "D3" = "2.416.66"
    "E3" = "  +1.29%  "
    "D4" = "1.00"
    "E4" = "  -0.09%  "
    "D5" = "562.22"
    "E5" = "  +2.08%  "
    "D6" = "142.86"
    "E6" = "  +0.93%  "
    "E7" = "  +0.07%  "
    "E8" = "  +1.64%  "
    "D9" = "2.412.17"
    "E9" = "  +1.02%  "
    "E10" = "  +1.32%  "
    "E11" = "  -2.16%  "
    "E12" = "  -0.26%  "
    "E13" = "  +0.21%  "
    "D14" = "25.65"
    "E14" = "  -0.76%  "
    "E15" = "  +0.72%  "
    "D16" = "2.852.46"
    "E16" = "  +1.26%  "
    "D17" = "62.113.76"
    "E17" = "  +1.07%  "
    "D18" = "2.413.77"
    "E18" = "  +1.13%  "
    "D19" = "11.32"
    "E19" = "  +2.15%  "
    "D20" = "4.17"
    "E20" = "  +0.38%  "
    "D21" = "323.24"
    "E21" = "  +0.28%  "
    "E22" = "  +3.08%  "
    "E23" = "  -0.07%  "
    "D24" = "65.77"
    "E24" = "  +2.08%  "
    "E25" = "  -1.52%  "
    "D26" = "9.02"
    "E26" = "  -1.91%  "
    "D27" = "577.52"
    "E27" = "  +4.72%  "
    "E28" = "  +1.20%  "
    "E29" = "  +0.16%  "
    "D30" = "0.0₃0947"
    "E30" = "  +3.53%  "
    "E31" = "  -0.28%  "
    "E32" = "  +1.29%  "
    "E34" = "  +0.92%  "
    "E35" = "  +0.39%  "
    "E36" = "  +0.14%  "
    "E37" = "  -2.90%  "
    "D40" = "152.07"
    "E40" = "  +4.32%  "
    "D41" = "18.68"
    "E41" = "  +0.69%  "
    "E42" = "  -6.88%  "
    "E43" = "  -0.14%  "
    "D44" = "2.29"
    "E44" = "  +1.99%  "
    "D45" = "148.53"
    "E45" = "  +0.06%  "
    "D47" = "0.0534"
    "E47" = "  +1.18%  "
    "D48" = "20.05"
    "E48" = "  -0.22%  "
    "D49" = "0.594"
    "E49" = "  +1.61%  "
    "D50" = "0.0919"
    "E50" = "  +1.37%  "
    "E51" = "  +1.96%  "
}

foreach ($ref in $updates.Keys) {
    $value = $updates[$ref]
    $cell = $ws.Range($ref)
    if ($forceText -contains $ref) {
        $cell.NumberFormat = "@"
        $cell.Value = $value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $value
    }
}
